$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Tables on slides 14, 15, 16: switch the applied table style from
#    "Medium Style 2 - Accent 1" ({ED391824-7844-4E10-B758-3E3B8517379B})
#    to the new style ({C87C01BA-15C2-4C36-8039-C650504A1B9C}).
# ---------------------------------------------------------------------------
$newTableStyleId = "{C87C01BA-15C2-4C36-8039-C650504A1B9C}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Theme: swap the presentation's colour scheme from "Red Violet"
#    (Integral theme) to the standard "Office" palette.
# ---------------------------------------------------------------------------
$slideOne = $p.Slides.Item(1)
$themeColors = $slideOne.ThemeColorScheme

# Order exposed by ThemeColorScheme.Item(1..12):
#   dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
#   accent6, hlink, folHlink
$officePalette = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officePalette[$i - 1]
}
